# Re-ran resolve and classify+summarise steps after changes to mapping file.
# This updates the "Range Status" sheet (Species (perc.) column removed /
# Species (no.) zeroed out) and the "High Priority break-up" sheet (the
# "Trend Different" row is gone and Trend New / IUCN rows carry new totals).

$wb = $excel.ActiveWorkbook

# --- Sheet 2: "Range Status" ---------------------------------------------
$wsRange = $wb.Worksheets.Item("Range Status")

$wsRange.Range("B2").Value = 0
$wsRange.Range("C2").ClearContents()

$wsRange.Range("B3").Value = 0
$wsRange.Range("C3").ClearContents()

$wsRange.Range("B4").Value = 0
$wsRange.Range("C4").ClearContents()

$wsRange.Range("B5").Value = 0
$wsRange.Range("C5").ClearContents()

$wsRange.Range("B6").Value = 0
$wsRange.Range("C6").ClearContents()

$wsRange.Range("B7").Value = 0
$wsRange.Range("C7").ClearContents()

# --- Sheet 5: "High Priority break-up" ------------------------------------
$wsBreakup = $wb.Worksheets.Item("High Priority break-up")

# "Trend Different" row disappears entirely; delete the whole row so
# "IUCN" shifts up from row 4 to row 3.
$wsBreakup.Rows.Item(3).Delete()

# New totals for the remaining two rows.
$wsBreakup.Range("B2").Value = 7
$wsBreakup.Range("C2").Value = 36.8
$wsBreakup.Range("D2").Value = 7
$wsBreakup.Range("E2").Value = 36.8

$wsBreakup.Range("B3").Value = 12
$wsBreakup.Range("C3").Value = 63.2
$wsBreakup.Range("D3").Value = 12
$wsBreakup.Range("E3").Value = 63.2
